## Team05Report.xlsx - "Add files via upload"
## Content changes:
##   1) Backlog sheet: swap the Owner initials for the JC04 / JC05 stories
##      (JC04 was "ik", JC05 was "tk" -> JC04 becomes "tk", JC05 becomes "ik").
##   2) Backlog sheet: a new backlog entry for story JC06 ("Display Operations")
##      is added, owned by "ik", status "Completed" - inserted as the new
##      row 4 (pushing the "No Backlogs" footer row down to row 5, and the
##      merged footer cell range down to A5:E5 accordingly).
##   3) The workbook's active sheet moves from "Team" to "Backlog", with the
##      selection on Backlog landing on E14.

$wb = $excel.ActiveWorkbook

$team = $wb.Worksheets.Item("Team")
$backlog = $wb.Worksheets.Item("Backlog")

## --- 1) swap the Owner initials on the two existing backlog rows ---------
$backlog.Range("D2").Value = "tk"
$backlog.Range("D3").Value = "ik"

## --- 2) insert the new JC06 / "Display Operations" backlog row ----------
$backlog.Rows.Item(4).Insert()

# Pick up the number formats / borders / fill from the row above (row 2),
# the same formatting every other backlog data row already uses.
$backlog.Range("A2:E2").Copy()
$backlog.Range("A4:E4").PasteSpecial(-4122)

$backlog.Range("A4").Value = 2
$backlog.Range("B4").Value = "JC06"
$backlog.Range("C4").Value = "Display Operations"
$backlog.Range("D4").Value = "ik"
$backlog.Range("E4").Value = "Completed"

## --- 3) move the active sheet/selection from Team to Backlog ------------
[void]$backlog.Activate()
[void]$backlog.Range("E14").Select()
